$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename header strings: "<Name>_old" -> "<Name>_FV2310", "<Name>_new" -> "<Name>_FV2404"
#    Columns A:J carry the "_old" headers, columns L:U carry the "_new" headers.
#    Column K ("diff") is untouched.
# ---------------------------------------------------------------------------
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $colOld = $i + 1       # columns A..J
    $colNew = $i + 12      # columns L..U
    $ws.Cells.Item(1, $colOld).Value = "$($baseNames[$i])_FV2310"
    $ws.Cells.Item(1, $colNew).Value = "$($baseNames[$i])_FV2404"
}

# ---------------------------------------------------------------------------
# 2. Turn the used range A1:U62 into an Excel Table ("Table1").
# ---------------------------------------------------------------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U62"), [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"
# Keep the table on the workbook's declared default style (same style the
# table would use if no explicit style name were stored at all).
$tbl.TableStyle = "TableStyleMedium9"

# ---------------------------------------------------------------------------
# 3. Freeze the header row (split after row 1, frozen state).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
